$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of key (column A value) -> new hash (column B value)
$updates = @{
    "01-010063TM" = "bf44ec96ab80cb1716583fd2713be6b3"
    "01-010063TC" = "7cd71806c3817a2788b411cc5dc0d07f"
    "01-010065TM" = "1c6b965ee60990ab7717581b4a83445d"
    "05-050203TC" = "764f754ffeb9b0002ad27d3f9529e20b"
    "03-030077TC" = "f1173a4f1a4db9a6ccec48b3aff58281"
    "03-030077TM" = "c217a913187da46d98037b67eb3ebdaf"
    "03-030077TP" = "0356828e5f648d949da58190375d702e"
    "05-050005TP" = "fa67257d9e82773e7b9d6f5b58515c14"
    "05-050001TC" = "edf7d3b415779cdb385a318df41f3e3e"
    "05-0709-070905BTC" = "afba4ee92bb44bede48ddf483ac24705"
    "01-010064TM" = "ece6eb734faed0dd6d9b51a279f5053d"
    "05-050201A" = "85a8a1badfc069c58b6cb6d4eac62f1a"
    "01-010064TC" = "a129a870088d76f781fe1f5950d3a8ba"
    "01-010064A" = "03f38022c575245c28fc04992de3c384"
    "05-050005A" = "e375d004872e7eac94fce210d9414135"
    "03-030078TC" = "a9ea093c40eaf3e1f00e4a1907276733"
    "03-030078TM" = "7f7ab1f8dc3ebc7cf76fcb6d6f79cd33"
    "03-030078TP" = "7196dec3ea8c8be2c644d2ff1202802c"
    "05-050208A" = "2567acf91628643a6f1b304994d0cb9f"
    "03-030076TC" = "e87098459d2c2631eb9417c08723ab77"
    "03-030076TM" = "40dab918aa390997041c69a02ac2fa13"
    "03-030076TP" = "aa3accc0757fd3e651b7e6d6f21f86c1"
    "01-010063A" = "ee5f9b6f034b61262ef8922f4d4f5ebd"
    "03-030032A" = "d878f735a89572d2273c1e98708e28dd"
    "03-030078A" = "012a3efc3a13ac4e2a1886c163d35e01"
    "03-030077A" = "3f27e8aadd43ec8a51d6e3542f7dce0e"
    "03-030076A" = "2e5ba69f5315bf502833e42bab0a83bb"
}

$colA = $ws.Range("A1:A962")
foreach ($key in $updates.Keys) {
    $cell = $colA.Find($key)
    if ($cell -ne $null) {
        $ws.Cells.Item($cell.Row, 2).Value2 = $updates[$key]
    }
}
